$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBody = "Beste klant,`nBedankt voor uw e-mail. Om uw vraag beter te kunnen begrijpen en u verder te kunnen helpen, zou ik graag wat meer informatie ontvangen. Kunt u mij uw ordernummer of referentienummer doorgeven, zodat ik kan controleren wat de status is van uw terugbetaling?`nIk kijk uit naar uw reactie.`nMet vriendelijke groet,`n[Naam] E-mailassistent `n[Bedrijfsnaam]"

$ws.Range("A25").Value = "Testmail #12: Ik heb nog geen geld terug."
$ws.Range("B25").Value = $newBody
$ws.Range("C25").Value = "Ik heb nog geen geld terug."
$ws.Range("D25").Value = "mailmind.test@zohomail.eu"
$ws.Range("E25").Value = "Retour / Terugbetaling"
$ws.Range("F25").Value = "2025-07-29 21:53:18"
$ws.Range("G25").Value = "Ja"
$ws.Range("H25").Value = "Nee"
$ws.Range("I25").Value = "Ja"
$ws.Range("J25").Value = "Nee"

$ws.Rows.Item(25).AutoFit()
